$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -19.9546393584243
$ws.Cells.Item(2, 3).Value = 2.364490549797686
$ws.Cells.Item(2, 4).Value = -19.9546393584243
$ws.Cells.Item(2, 5).Value = -19.9546393584243
$ws.Cells.Item(2, 6).Value = -19.9546393584243
$ws.Cells.Item(2, 7).Value = -19.9546393584243
$ws.Cells.Item(2, 8).Value = -19.9546393584243
$ws.Cells.Item(2, 9).Value = -19.9546393584243
$ws.Cells.Item(2, 10).Value = -19.9546393584243
$ws.Cells.Item(2, 11).Value = -19.9546393584243

$ws.Cells.Item(3, 2).Value = -19.9546393584243
$ws.Cells.Item(3, 3).Value = -19.9546393584243
$ws.Cells.Item(3, 4).Value = -19.9546393584243
$ws.Cells.Item(3, 5).Value = -19.9546393584243
$ws.Cells.Item(3, 6).Value = -19.9546393584243
$ws.Cells.Item(3, 7).Value = -19.9546393584243
$ws.Cells.Item(3, 8).Value = -19.9546393584243
$ws.Cells.Item(3, 9).Value = 2.313818403471394
$ws.Cells.Item(3, 10).Value = -19.9546393584243
$ws.Cells.Item(3, 11).Value = -19.9546393584243

$ws.Cells.Item(4, 2).Value = -19.9546393584243
$ws.Cells.Item(4, 3).Value = 2.218495164842582
$ws.Cells.Item(4, 4).Value = 4.321926746069378
$ws.Cells.Item(4, 5).Value = -19.9546393584243
$ws.Cells.Item(4, 6).Value = 2.451542027005565
$ws.Cells.Item(4, 7).Value = -19.9546393584243
$ws.Cells.Item(4, 8).Value = 1.786928435219915
$ws.Cells.Item(4, 9).Value = -19.9546393584243
$ws.Cells.Item(4, 10).Value = 2.3465690588064
$ws.Cells.Item(4, 11).Value = -19.9546393584243

$ws.Cells.Item(5, 2).Value = -19.9546393584243
$ws.Cells.Item(5, 3).Value = 0.8060976416894708
$ws.Cells.Item(5, 4).Value = -19.9546393584243
$ws.Cells.Item(5, 5).Value = -19.9546393584243
$ws.Cells.Item(5, 6).Value = -19.9546393584243
$ws.Cells.Item(5, 7).Value = 2.128548099568577
$ws.Cells.Item(5, 8).Value = -19.9546393584243
$ws.Cells.Item(5, 9).Value = -19.9546393584243
$ws.Cells.Item(5, 10).Value = -19.9546393584243
$ws.Cells.Item(5, 11).Value = -19.9546393584243

$ws.Cells.Item(6, 2).Value = -19.9546393584243
$ws.Cells.Item(6, 3).Value = -19.9546393584243
$ws.Cells.Item(6, 4).Value = -19.9546393584243
$ws.Cells.Item(6, 5).Value = -19.9546393584243
$ws.Cells.Item(6, 6).Value = -19.9546393584243
$ws.Cells.Item(6, 7).Value = -19.9546393584243
$ws.Cells.Item(6, 8).Value = -19.9546393584243
$ws.Cells.Item(6, 9).Value = -19.9546393584243
$ws.Cells.Item(6, 10).Value = -19.9546393584243
$ws.Cells.Item(6, 11).Value = -19.9546393584243

$ws.Cells.Item(7, 2).Value = 3.084593816150727
$ws.Cells.Item(7, 3).Value = -19.9546393584243
$ws.Cells.Item(7, 4).Value = -19.9546393584243
$ws.Cells.Item(7, 5).Value = -19.9546393584243
$ws.Cells.Item(7, 6).Value = -19.9546393584243
$ws.Cells.Item(7, 7).Value = -19.9546393584243
$ws.Cells.Item(7, 8).Value = -19.9546393584243
$ws.Cells.Item(7, 9).Value = -19.9546393584243
$ws.Cells.Item(7, 10).Value = -19.9546393584243
$ws.Cells.Item(7, 11).Value = -19.9546393584243

$ws.Cells.Item(8, 2).Value = -19.9546393584243
$ws.Cells.Item(8, 3).Value = -19.9546393584243
$ws.Cells.Item(8, 4).Value = -19.9546393584243
$ws.Cells.Item(8, 5).Value = 2.909863636479878
$ws.Cells.Item(8, 6).Value = -19.9546393584243
$ws.Cells.Item(8, 7).Value = -19.9546393584243
$ws.Cells.Item(8, 8).Value = -19.9546393584243
$ws.Cells.Item(8, 9).Value = -19.9546393584243
$ws.Cells.Item(8, 10).Value = -19.9546393584243
$ws.Cells.Item(8, 11).Value = -19.9546393584243

$ws.Cells.Item(9, 2).Value = 3.525676659734234
$ws.Cells.Item(9, 3).Value = -19.9546393584243
$ws.Cells.Item(9, 4).Value = -19.9546393584243
$ws.Cells.Item(9, 5).Value = -19.9546393584243
$ws.Cells.Item(9, 6).Value = -19.9546393584243
$ws.Cells.Item(9, 7).Value = -19.9546393584243
$ws.Cells.Item(9, 8).Value = -19.9546393584243
$ws.Cells.Item(9, 9).Value = -19.9546393584243
$ws.Cells.Item(9, 10).Value = -19.9546393584243
$ws.Cells.Item(9, 11).Value = -19.9546393584243

$ws.Cells.Item(10, 2).Value = -19.9546393584243
$ws.Cells.Item(10, 3).Value = -19.9546393584243
$ws.Cells.Item(10, 4).Value = -19.9546393584243
$ws.Cells.Item(10, 5).Value = -19.9546393584243
$ws.Cells.Item(10, 6).Value = -19.9546393584243
$ws.Cells.Item(10, 7).Value = -19.9546393584243
$ws.Cells.Item(10, 8).Value = -19.9546393584243
$ws.Cells.Item(10, 9).Value = 1.519623594447283
$ws.Cells.Item(10, 10).Value = -19.9546393584243
$ws.Cells.Item(10, 11).Value = 2.21042495435323

$ws.Cells.Item(11, 2).Value = -19.9546393584243
$ws.Cells.Item(11, 3).Value = -19.9546393584243
$ws.Cells.Item(11, 4).Value = -19.9546393584243
$ws.Cells.Item(11, 5).Value = 1.838789715223383
$ws.Cells.Item(11, 6).Value = -19.9546393584243
$ws.Cells.Item(11, 7).Value = 2.601357921606791
$ws.Cells.Item(11, 8).Value = -19.9546393584243
$ws.Cells.Item(11, 9).Value = -19.9546393584243
$ws.Cells.Item(11, 10).Value = -19.9546393584243
$ws.Cells.Item(11, 11).Value = 1.386402157676419

$ws.Cells.Item(12, 2).Value = -19.9546393584243
$ws.Cells.Item(12, 3).Value = -19.9546393584243
$ws.Cells.Item(12, 4).Value = -19.9546393584243
$ws.Cells.Item(12, 5).Value = -19.9546393584243
$ws.Cells.Item(12, 6).Value = -19.9546393584243
$ws.Cells.Item(12, 7).Value = -19.9546393584243
$ws.Cells.Item(12, 8).Value = -19.9546393584243
$ws.Cells.Item(12, 9).Value = -19.9546393584243
$ws.Cells.Item(12, 10).Value = -19.9546393584243
$ws.Cells.Item(12, 11).Value = -19.9546393584243

$ws.Cells.Item(13, 2).Value = -19.9546393584243
$ws.Cells.Item(13, 3).Value = -19.9546393584243
$ws.Cells.Item(13, 4).Value = -19.9546393584243
$ws.Cells.Item(13, 5).Value = 1.705346520286392
$ws.Cells.Item(13, 6).Value = -19.9546393584243
$ws.Cells.Item(13, 7).Value = -19.9546393584243
$ws.Cells.Item(13, 8).Value = -19.9546393584243
$ws.Cells.Item(13, 9).Value = -19.9546393584243
$ws.Cells.Item(13, 10).Value = 2.277902483928688
$ws.Cells.Item(13, 11).Value = 1.624216547323618

$ws.Cells.Item(14, 2).Value = -19.9546393584243
$ws.Cells.Item(14, 3).Value = -19.9546393584243
$ws.Cells.Item(14, 4).Value = -19.9546393584243
$ws.Cells.Item(14, 5).Value = -19.9546393584243
$ws.Cells.Item(14, 6).Value = -19.9546393584243
$ws.Cells.Item(14, 7).Value = -19.9546393584243
$ws.Cells.Item(14, 8).Value = -19.9546393584243
$ws.Cells.Item(14, 9).Value = -19.9546393584243
$ws.Cells.Item(14, 10).Value = -19.9546393584243
$ws.Cells.Item(14, 11).Value = 2.084450039262292

$ws.Cells.Item(15, 2).Value = -19.9546393584243
$ws.Cells.Item(15, 3).Value = -19.9546393584243
$ws.Cells.Item(15, 4).Value = -19.9546393584243
$ws.Cells.Item(15, 5).Value = -19.9546393584243
$ws.Cells.Item(15, 6).Value = -19.9546393584243
$ws.Cells.Item(15, 7).Value = -19.9546393584243
$ws.Cells.Item(15, 8).Value = -19.9546393584243
$ws.Cells.Item(15, 9).Value = -19.9546393584243
$ws.Cells.Item(15, 10).Value = -19.9546393584243
$ws.Cells.Item(15, 11).Value = -19.9546393584243

$ws.Cells.Item(16, 2).Value = -19.9546393584243
$ws.Cells.Item(16, 3).Value = -19.9546393584243
$ws.Cells.Item(16, 4).Value = -19.9546393584243
$ws.Cells.Item(16, 5).Value = -19.9546393584243
$ws.Cells.Item(16, 6).Value = -19.9546393584243
$ws.Cells.Item(16, 7).Value = -19.9546393584243
$ws.Cells.Item(16, 8).Value = -19.9546393584243
$ws.Cells.Item(16, 9).Value = -19.9546393584243
$ws.Cells.Item(16, 10).Value = 2.321098972924648
$ws.Cells.Item(16, 11).Value = -19.9546393584243

$ws.Cells.Item(17, 2).Value = -19.9546393584243
$ws.Cells.Item(17, 3).Value = -0.3623411028510972
$ws.Cells.Item(17, 4).Value = -19.9546393584243
$ws.Cells.Item(17, 5).Value = -19.9546393584243
$ws.Cells.Item(17, 6).Value = -19.9546393584243
$ws.Cells.Item(17, 7).Value = -19.9546393584243
$ws.Cells.Item(17, 8).Value = 0.6080064788684003
$ws.Cells.Item(17, 9).Value = 0.8775591801450574
$ws.Cells.Item(17, 10).Value = 1.29968863564255
$ws.Cells.Item(17, 11).Value = -19.9546393584243

$ws.Cells.Item(18, 2).Value = -19.9546393584243
$ws.Cells.Item(18, 3).Value = -19.9546393584243
$ws.Cells.Item(18, 4).Value = -19.9546393584243
$ws.Cells.Item(18, 5).Value = -19.9546393584243
$ws.Cells.Item(18, 6).Value = -19.9546393584243
$ws.Cells.Item(18, 7).Value = -19.9546393584243
$ws.Cells.Item(18, 8).Value = 0.510114767739541
$ws.Cells.Item(18, 9).Value = 0.8844513763153681
$ws.Cells.Item(18, 10).Value = 1.381418752041407
$ws.Cells.Item(18, 11).Value = -19.9546393584243

$ws.Cells.Item(19, 2).Value = -19.9546393584243
$ws.Cells.Item(19, 3).Value = -19.9546393584243
$ws.Cells.Item(19, 4).Value = -19.9546393584243
$ws.Cells.Item(19, 5).Value = -19.9546393584243
$ws.Cells.Item(19, 6).Value = -19.9546393584243
$ws.Cells.Item(19, 7).Value = -19.9546393584243
$ws.Cells.Item(19, 8).Value = 1.939478603981972
$ws.Cells.Item(19, 9).Value = 2.164716273970559
$ws.Cells.Item(19, 10).Value = -19.9546393584243
$ws.Cells.Item(19, 11).Value = -19.9546393584243

$ws.Cells.Item(20, 2).Value = -19.9546393584243
$ws.Cells.Item(20, 3).Value = 1.935549573367839
$ws.Cells.Item(20, 4).Value = -19.9546393584243
$ws.Cells.Item(20, 5).Value = -19.9546393584243
$ws.Cells.Item(20, 6).Value = 3.860960589483754
$ws.Cells.Item(20, 7).Value = -19.9546393584243
$ws.Cells.Item(20, 8).Value = 2.189704278646093
$ws.Cells.Item(20, 9).Value = 1.997741093120474
$ws.Cells.Item(20, 10).Value = -19.9546393584243
$ws.Cells.Item(20, 11).Value = 2.44194269024849

$ws.Cells.Item(21, 2).Value = -19.9546393584243
$ws.Cells.Item(21, 3).Value = 1.942828352598168
$ws.Cells.Item(21, 4).Value = -19.9546393584243
$ws.Cells.Item(21, 5).Value = 2.497316507436802
$ws.Cells.Item(21, 6).Value = -19.9546393584243
$ws.Cells.Item(21, 7).Value = 3.256804731393366
$ws.Cells.Item(21, 8).Value = 2.379367895418214
$ws.Cells.Item(21, 9).Value = -19.9546393584243
$ws.Cells.Item(21, 10).Value = -19.9546393584243
$ws.Cells.Item(21, 11).Value = -19.9546393584243

Write-Host "Updated PSSM values for supplemental figures"